$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 374, shifting the existing rows 374:399 down to 375:400
$ws.Rows.Item(374).Insert()

# Populate the newly inserted row 374 with the new record
$ws.Cells.Item(374, 1).Value = 10
$ws.Cells.Item(374, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(374, 3).Value = "La Araucanía"
$ws.Cells.Item(374, 4).Value = 44931
$ws.Cells.Item(374, 5).Value = 9
$ws.Cells.Item(374, 6).Value = 100114013
$ws.Cells.Item(374, 7).Value = "Zanahoria"
$ws.Cells.Item(374, 8).Value = "Sin especificar"
$ws.Cells.Item(374, 9).Value = "Primera"
$ws.Cells.Item(374, 10).Value = 50
$ws.Cells.Item(374, 11).Value = 14000
$ws.Cells.Item(374, 12).Value = 15000
$ws.Cells.Item(374, 13).Value = 14600
$ws.Cells.Item(374, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(374, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(374, 16).Value = 584
$ws.Cells.Item(374, 17).Value = 25
$ws.Cells.Item(374, 18).Value = "Hortaliza"
